$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    23847147.04472806,23874991.41112009,23903205.79653071,23931794.92695915,23960763.57587729,23990116.56381527,24019858.75788571,24049995.07124395,24080530.46248226,24111469.93495544,
    24142818.53603546,24174581.3562934,24206763.52860603,24239370.22718549,24272406.66653002,24305878.100294,24339789.82007609,24374147.15412356,24408955.46595238,24444220.15288147,
    24479946.64448103,24516140.40093416,24552806.91131188,24589951.6917619,24627580.28361125,24665698.25138437,24704311.18073734,24743424.67631048,24783044.35950117,24823175.86615981,
    24863824.84421166,24904996.95120861,24946697.85181468,24988933.2152301,25031708.71255923,25075030.01412812,25118902.78675813,25163332.69100262,25208325.37835471,25253886.48843398,
    25300021.64616132,25346736.45893168,25394036.51379494,25441927.37465567,25490414.57950387,25539503.63768839,25589200.02724683,25639509.19230458,25690436.54055808,25741987.44085667,
    25794167.2208988,25846981.16505835,25900434.51235803,25954532.45460674,26009280.13471844,26064682.64523081,26120745.02704196,26177472.26838401,26234869.30405264,26292941.0149118,
    26351692.22769353,26411127.71511163,26471252.19630997,26532070.33766384,26593586.75395503,26655806.00993897,26718732.62232358,26782371.06217851,26846725.75779273,26911801.09799841,
    26977601.43597822,27044131.09357226,27111394.36610027,27179395.52771394,27248138.83729247,27317628.54489437,27387868.89877667,27458864.15299129,27530618.57556751,27603136.4572873,
    27676422.1210595,27750479.93189586,27825314.30749211,27900929.72941316,27977330.75488183,28054522.02916638,28132508.2985618,28211294.42395726,28290885.39497868,28371286.3446953,
    28452502.56487534,28534539.52177386,28617402.87243428,28701098.48148199,28785632.43838617,28871011.07516477,28957240.98450333,29044329.0382583,29132282.40631168,29221108.57574208,
    29310815.37027584,29401410.96997913,29492903.9311498,29585303.20636715,29678618.16465458,29772858.61170942,29868034.8101523,29964157.49974759,30061237.9175443,30159287.81788684,
    30258319.49224365,30358345.7888005,30459380.13176602,30561436.54033527,30664529.6472579,30768674.71695742,30873887.66314789,30980185.06589518,31087584.18807061,31196102.99114579,
    31305760.15027792,31416575.0686371,31528567.89092808,31641759.51606056,31756171.60892494,31871826.61123125,31988747.75137269,32106959.05327704,32226485.34421189,32347352.26151287,
    32469586.25820727,32593214.60750816,32718265.40615817,32844767.57660513,32972750.86799589,33102245.85597766,33233283.94130185,33365897.34722649,33500119.11572067,33635983.10247591,
    33773523.97073491,33912777.18395186,34053778.99730267,34196566.44806799,34341177.34491529,34487650.25611113,34636024.49669856,34786340.11467792,34938637.87623462,35092959.25005923,
    35249346.39081087,35407842.12177686,35568489.91678584,35731333.88143389,35896418.73368751,36063789.78392825,36233492.914508,36405574.55888569,36580081.68041775,36757061.75087781,
    36936562.72878141,37118633.037594,37303321.54390127,37490677.53562117,37680750.70033896,37873591.10384568,38069249.16896141,38267775.65472377,38469221.63602312,38673638.4837633,
    38881077.84562747,39091591.62752698,39305231.97580935,39522051.26030051,39742102.05825403,39965437.13927872,40192109.45131287,40422172.10771178,40655678.37551197,40892681.66493338,
    41133235.52017681,41377393.61157268,41625209.72913142,41876737.77754532,42132031.77268612,42391145.83963969,42654134.21231746,42921051.23467675,43191951.36358361,43466889.17334348,
    43745919.3619244,44029096.75889341,44316476.33508156,44608113.21399226,44904062.68496133,45204380.21807632,45509121.48085767,45818342.35670169,46132098.96508235,46450447.68350643,
    46773445.1712122,47101148.39460094,47433614.65437743,47770901.61441064,48113067.33226924,48460170.29142076,48812269.43507066,49169424.20161809,49531694.56170175,49899141.05680934,
    50271824.83942167,50649807.71466286,51033152.18342477,51421921.48693696,51816179.65274837,52215991.54209024,52621422.89858755,53032540.39828666,53449411.70096651,53872105.50270241,
    54300691.58964821,54735240.89300796,55175825.54516312,55622518.93692695,56075395.77589414,56534532.1458571,57000005.56725975,57471895.05866037,57950281.19917607,58435246.19188239,
    58926873.9281408,59425250.05282915,59930462.03045028,60442599.21209312,60961752.90322308,61488016.43227731,62021485.22004124,62562256.84978286,63110431.13812105,63666110.20660347,
    64229398.55397034,64800403.12907887,65379233.40446191,65966001.45049495,66560822.01014215,67163812.57425304,67775093.45737693,68394787.8740624,69023022.0156059,69659925.12721008,
    70305629.58551115,70960270.97642945,71623988.17329562,72296923.41519858,72979222.38549936,73671034.29044865,74372511.93784094,75083811.81563337,75805094.17045115,76536523.08589409,
    77278266.56055468,78030496.58564715,78793389.22214416,79567124.67730515,80351887.38047567,81147866.05802564,81955253.80728799,82774248.1693459,83605051.20051104,84447869.54232125,
    85302914.48987767,86170402.05832821,87050553.04729517,87943593.10303064,88849752.77807102,89769267.58815263,90702378.06613322,91649329.81265272,92610373.54325642,93585765.13168265,
    94575765.64901222,95580641.39835486,96600663.94474001,97636110.13986093,98687262.14130944,99754407.42592299,100837838.7968535,101937854.3839507,103054757.6370423,104188857.3116781,
    105340467.4468914,106509907.3345207,107697501.4796216,108903579.5514862,110128476.3247782,111372531.6102823,112636090.1747555,113919501.6493623,115223120.4261658,116547305.542147,
    117892420.5502123,118493488.3365386,119105245.5182933,119727895.7464864,120361644.734226,121006700.1313354,121663271.3904382,122331569.6242577,123011807.4538838,123704198.8477671,
    124408958.9512142,125126303.9061668,125856450.6610604,126599616.7705767,127356020.1851195,128125879.0298593,128909411.3732213,129706834.9847049,130518367.0819537,131344224.067022,
    132184621.25181,133039772.5726727,133909890.2942418,134795184.7025359,135695863.7874664,136612132.9148958,137544194.488437,138492247.6012331,139456487.6779971,140437106.1076427,
    141434289.8668841,142448221.1352315,143479076.9018635,144527028.5649103,145592241.5237325,146674874.764838,147775080.4421349,148893003.452269,150028781.0058571,151182542.1954771,
    152354407.5613317,153544488.6555575,154752887.6062018,155979696.6819434,157224997.8586753,158488862.3891238,159771350.3767124,161072510.3549247,162392378.8734553,163730980.0924714,
    165088325.3863359,166464412.9581626,167859227.4666012,169272739.6662508,170704906.0631195,172155668.586537,171506168.1963763,170850473.3456461,170188505.7449939,169520185.9414427,
    168845433.3193823,168164166.1009226,167476301.3455604,166781754.9491317,166080441.6420246,165372274.9866598,164657167.3742555,163935030.0209219,163205772.9631557,162469305.0528313,
    161725533.9518209,160974366.1264063,160215706.8416887,159449460.1562346,158675528.9172517,157893814.7566255,157104218.0882039,156306638.1067672,155500972.7891782,154687118.8982647,
    153864971.9900482,153034426.4249938,152195375.3840225,151347710.8900925,150491323.8362196,149626104.0208774,148751940.1917795,147868720.0991157,146976330.5593709,146074657.5309219,
    145163586.2026577,144243001.0969252,143312786.1881467,142372825.0384949,141423000.952043,140463197.1488325,139493296.9603097,138513184.0475898,137522742.6439921,136521857.8232664,
    135510415.7948952,134488304.2277967,133455412.6036869,132411632.601266,131356858.5122906,130290987.6904624,129213921.033919,128125563.5019443,127025824.6663294,125914619.2976032,
    124791867.9861233,123657497.7977705,122511442.9637158,121353645.6034454,120184056.4799189,119002635.7854164,117809353.9562926,116604192.5145049,115387144.9334272,114158217.5250874,
    112917430.3456007,111664818.1151924,110400431.1488349,109124336.2931532,107836617.8649002,106537378.5859524,105226740.5094532,103904845.9314228,102571858.2818735,101227962.9892166,
    99873368.31153184,98508306.12809022,97133032.68438391,95747829.28382841,94353002.9192552,92948886.83733094,91535841.02909866,90114252.63996556,88684536.29264233,87247134.31678464,
    85802516.87939411,84351182.01040336,82893655.51830085,81430490.79114088,79962268.47883373,78489596.05321638,77013107.24306276,75533461.34190074,74051342.38725519,72567458.21073027,
    71082539.35916865
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
